$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing data rows (2-11) down to (3-12)
$ws.Range("A2:C2").Insert()

# Force the new cells to Text format first so values aren't reinterpreted as numbers/dates
$ws.Range("A2:C2").NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "+553187596479"
$ws.Cells.Item(2, 2).Value = "31"
$ws.Cells.Item(2, 3).Value = "2024-10-22"

# Now match the visual formatting (font/fill/border/alignment/number format) of the data rows
$ws.Range("A3:C3").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
